$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values are being regenerated (Strike# -> K).
# Update the G column cells for rows 2-10 per the new computed values.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 1
$ws.Range("G10").Value = 2
